# Update the "Förändrad" (C) column date values from 45183 (2023-09-14)
# to 45184 (2023-09-15) for rows 2 through 16 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
